# "Fin de la Phase 1" — populate estimated-hours (E) for the remaining
# tasks, restore the shared IF() formulas in column H, add the thread of
# review comments left on column C, and leave the sheet scrolled/selected
# where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plannification globale")

# --- Column E: fill in the estimated-hours values that were missing ---
$ws.Range("E40").Value = 34.35
$ws.Range("E41").Value = 32
$ws.Range("E42").Value = 33
$ws.Range("E45").Value = 2.3
$ws.Range("E46").Value = "1,4,8,11,18,34"
$ws.Range("E47").Value = 21
$ws.Range("E48").Value = 21
$ws.Range("E49").Value = 21
$ws.Range("E50").Value = 32.33
$ws.Range("E51").Value = "34, 18, 15, 11, 7, 4, 1"
$ws.Range("E52").Value = "12, 35"
$ws.Range("E53").Value = 36

# --- Column H: re-enter the difficulty formula across each contiguous
#     block so Excel collapses it back into a shared formula group ---
$ws.Range("H9:H14").Formula = '=IF(G9="facile",1,IF(G9="moyen",5,IF(G9="difficile",10,7.5)))'
$ws.Range("H16:H53").Formula = '=IF(G16="facile",1,IF(G16="moyen",5,IF(G16="difficile",10,7.5)))'

# --- Threaded review comments added on column C ---
$ws.Range("C9").AddCommentThreaded("Affichage incluant toujours : Incluant CSS, JS et html") | Out-Null
$ws.Range("C10").AddCommentThreaded("Création des setters pour la BD") | Out-Null
$ws.Range("C11").AddCommentThreaded("gestion de l'authentification dans cette étape") | Out-Null
$ws.Range("C14").AddCommentThreaded("gestion des formulaires et affichage des erreurs le cas échéant") | Out-Null
$ws.Range("C15").AddCommentThreaded("écriture dans la BD") | Out-Null
$ws.Range("C17").AddCommentThreaded("doit lire dans BD évènement pour afficher les évènements en lien avec le userID") | Out-Null
$ws.Range("C20").AddCommentThreaded("et création d'une disponibilité pour le créateur") | Out-Null
$ws.Range("C24").AddCommentThreaded("Redirection vers acceuil si client non authentifié essaie d'accéder à une page") | Out-Null
$ws.Range("C28").AddCommentThreaded("Initialisation des paramètres : nb de genérations maximales, case horaire de l'évènement et les disponibilités des participants") | Out-Null
$ws.Range("C30").AddCommentThreaded("Ne doit pas avoir de lien vers le problème directement.") | Out-Null
$ws.Range("C38").AddCommentThreaded("Un seul choix pour la disponibilité. Réserve automatiquement la plage horaire. La plage devient indisponible pour les autres usagers.") | Out-Null
$ws.Range("C39").AddCommentThreaded("Possibilité d'ajouter une préférence. Doit être tenu en compte comme une variable par le ScheduleOptimizer") | Out-Null
$ws.Range("C44").AddCommentThreaded("Envoyer un JSON avec les données à être affichée pour l'affichage. ET pour le calcul du script algoGen") | Out-Null
$ws.Range("C52").AddCommentThreaded("Fonctionnalité similaire à Doodle. Un sondage seulement.") | Out-Null

# --- Leave the view where the author ended up ---
$ws.Activate()
$ws.Range("B38").Select()
